# Re-run of the experiment task-order generator: each task-order sheet is
# regenerated (new timestamp-suffixed stim filenames) and the tabs come back
# out in a new shuffled order. We reorder tabs to match, rename them, and
# overwrite the task_order column with the freshly generated filenames.

$wb = $excel.ActiveWorkbook

# --- 1. Reorder the tabs -----------------------------------------------
# Original order:  1=GNG  2=NB  3=RS  4=TOL  5=vSAT
# Target order:     vSAT, TOL, NB, GNG, RS
# Move the vSAT sheet (currently last) to the front, then TOL in behind it,
# leaving NB/GNG/RS in their relative order after.
$wsVSAT = $wb.Worksheets.Item(5)
$wsVSAT.Move($wb.Worksheets.Item(1))

$wsTOL = $wb.Worksheets.Item("TOL_TO-16512555696484015")
$wsTOL.Move($wb.Worksheets.Item(2))

# Now order is: vSAT, TOL, GNG, NB, RS -- bring NB ahead of GNG.
$wsNB = $wb.Worksheets.Item("NB_TO-16512555695848374")
$wsNB.Move($wb.Worksheets.Item(3))

# Final order should now be: vSAT, TOL, NB, GNG, RS

# --- 2. Rename the tabs --------------------------------------------------
$wb.Worksheets.Item(1).Name = "vSAT_TO-1651588999806087"
$wb.Worksheets.Item(2).Name = "TOL_TO-16515889998502576"
$wb.Worksheets.Item(3).Name = "NB_TO-16515890024288406"
$wb.Worksheets.Item(4).Name = "GNG_TO-16515890024601283"
$wb.Worksheets.Item(5).Name = "RS_TO-16515890024601283"

# --- 3. Refresh the generated stim filenames in each task_order sheet ----

# vSAT sheet (position 1) -- SAT/vSAT stims
$ws = $wb.Worksheets.Item(1)
$ws.Range("B2").Value = "SAT_stims-16515889997592113.csv"
$ws.Range("B3").Value = "vSAT_stims-16515889997904606.csv"
$ws.Range("B4").Value = "vSAT_stims-1651588999774835.csv"
$ws.Range("B5").Value = "SAT_stims-1651588999743583.csv"

# TOL sheet (position 2) -- MM/ZM stims
$ws = $wb.Worksheets.Item(2)
$ws.Range("B2").Value = "MM_stims-16515889998217103.csv"
$ws.Range("B3").Value = "ZM_stims-1651588999806087.csv"
$ws.Range("B4").Value = "MM_stims-16515889998346305.csv"
$ws.Range("B5").Value = "ZM_stims-16515889998230846.csv"
$ws.Range("B6").Value = "MM_stims-16515889998502576.csv"
$ws.Range("B7").Value = "ZM_stims-16515889998346305.csv"

# NB sheet (position 3) -- ZB/OB/TB stims
$ws = $wb.Worksheets.Item(3)
$ws.Range("B2").Value = "TB-1651589002413251.csv"
$ws.Range("B3").Value = "OB-16515890015442252.csv"
$ws.Range("B4").Value = "TB-16515890019892068.csv"
$ws.Range("B5").Value = "ZB-match_7-1651589000303104.csv"
$ws.Range("B6").Value = "TB-16515890023350904.csv"
$ws.Range("B7").Value = "ZB-match_8-16515890003677714.csv"
$ws.Range("B8").Value = "ZB-match_6-16515890000610185.csv"
$ws.Range("B9").Value = "OB-16515890017963827.csv"
$ws.Range("B10").Value = "OB-16515890013057544.csv"

# GNG sheet (position 4) -- go/GNG stims
$ws = $wb.Worksheets.Item(4)
$ws.Range("B2").Value = "go_stims-16515890024288406.csv"
$ws.Range("B3").Value = "GNG_stims-1651589002444465.csv"
$ws.Range("B4").Value = "go_stims-1651589002444465.csv"
$ws.Range("B5").Value = "GNG_stims-16515890024601283.csv"

# RS sheet (position 5) -- unchanged (eyes open / eyes closed)
